$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 181, shifting existing rows 181-227 down to 182-228
$ws.Rows("181:181").Insert()

# Populate the newly inserted row 181 with the new data record
$ws.Range("A181").Value = 4
$ws.Range("B181").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C181").Value = "Los Lagos"
$ws.Range("D181").Value = 44841
$ws.Range("E181").Value = 10
$ws.Range("F181").Value = 100112009
$ws.Range("G181").Value = "Acelga"
$ws.Range("H181").Value = "Sin especificar"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 200
$ws.Range("K181").Value = 3000
$ws.Range("L181").Value = 3500
$ws.Range("M181").Value = 3250
$ws.Range("N181").Value = "$/docena de atados (4 kilos)"
$ws.Range("O181").Value = "Región del Maule"
$ws.Range("P181").Value = 812
$ws.Range("Q181").Value = 4
$ws.Range("R181").Value = "Hortaliza"
